$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "San Fernando"
$ws.Range("B1").Value = "Sevilla"
$ws.Range("C1").Value = "Sevilla"
$ws.Range("D1").Value = 41004

$ws.Range("H1:J2").NumberFormat = "@"

$ws.Range("H1").Value = "37.3134419"
$ws.Range("I1").Value = "-4.869306"
$ws.Range("J1").Value = "-4.869306"

$ws.Range("H2").Value = "43.2687376"
$ws.Range("I2").Value = "-2.9404136"
$ws.Range("J2").Value = "-2.9404136"

$ws.Range("H1:J2").NumberFormat = "General"

$ws.Range("D2").Select()
